$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.864.97'
$ws.Range("E2").Value = '  +0.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.751.79'
$ws.Range("E3").Value = '  +0.64%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9986'
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.39'
$ws.Range("E5").Value = '  +0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9988'
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5136'
$ws.Range("E7").Value = '  +4.86%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.42'
$ws.Range("E8").Value = '  -1.98%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2674'
$ws.Range("E9").Value = '  +3.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06187'
$ws.Range("E10").Value = '  +2.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.792.44'
$ws.Range("E11").Value = '  +2.87%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06946'
$ws.Range("E12").Value = '  +1.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.43'
$ws.Range("E13").Value = '  +4.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6281'
$ws.Range("E14").Value = '  +11.69%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.486'
$ws.Range("E15").Value = '  +0.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.96'
$ws.Range("E16").Value = '  +3.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9993'
$ws.Range("E17").Value = '  -0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9993'
$ws.Range("E18").Value = '  -0.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.859.88'
$ws.Range("E19").Value = '  +0.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.62'
$ws.Range("E20").Value = '  +1.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006664'
$ws.Range("E21").Value = '  +1.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.988.47'
$ws.Range("E22").Value = '  +1.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.060'
$ws.Range("E23").Value = '  +0.84%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.271'
$ws.Range("E24").Value = '  +4.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.172'
$ws.Range("E25").Value = '  +3.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '136.43'
$ws.Range("E26").Value = '  -0.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.479'
$ws.Range("E27").Value = '  +1.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.14'
$ws.Range("E28").Value = '  +3.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.783'
$ws.Range("E29").Value = '  -2.18%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.82'
$ws.Range("E30").Value = '  +1.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08269'
$ws.Range("E31").Value = '  +3.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.690'
$ws.Range("E32").Value = '  -0.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.408'
$ws.Range("E33").Value = '  +0.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04392'
$ws.Range("E34").Value = '  -0.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.642'
$ws.Range("E35").Value = '  +1.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").Value = '  +2.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6048'
$ws.Range("E37").Value = '  +2.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.653'
$ws.Range("E38").Value = '  -0.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01558'
$ws.Range("E39").Value = '  +3.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.946'
$ws.Range("E40").Value = '  +6.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9992'
$ws.Range("E41").Value = '  -0.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '102.19'
$ws.Range("E42").Value = '  -0.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3844'
$ws.Range("E43").Value = '  +3.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7469'
$ws.Range("E44").Value = '  +3.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.889'
$ws.Range("E45").Value = '  -4.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05485'
$ws.Range("E46").Value = '  +5.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1100'
$ws.Range("E47").Value = '  +1.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.973'
$ws.Range("E48").Value = '  +3.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.09'
$ws.Range("E49").Value = '  +0.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.71'
$ws.Range("E50").Value = '  +1.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").Value = '  +0.46%  '
